# The commit swaps the two theme parts in the deck:
#   ppt/theme/theme1.xml  (the theme actually used by the slide master /
#                           all slides) goes from the "Integral" / "Red
#                           Violet" colour scheme to the stock "Office
#                           Theme" colour scheme.
#   ppt/theme/theme2.xml  (the theme used by the notes master) goes the
#                           other way, from "Office Theme" to "Integral"
#                           / "Red Violet".
#
# The font scheme and format scheme (fills/lines/effects) are byte-for-
# byte identical between the two theme parts, so the only observable
# difference is the 12-slot colour scheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink). We reproduce that through ThemeColorScheme.Colors(i).RGB
# - the automation surface PowerPoint exposes for editing a theme's colour
# scheme - on the slide master, which is the theme part that is visibly in
# effect across every slide in the deck.

function New-BGR([int]$r, [int]$g, [int]$b) {
    # OLE_COLOR / VBA RGB() packs as 0x00BBGGRR.
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette for ppt/theme/theme1.xml ("Office Theme").
$officeTheme = @{
    1  = (New-BGR 0x00 0x00 0x00)  # dk1      000000
    2  = (New-BGR 0xFF 0xFF 0xFF)  # lt1      FFFFFF
    3  = (New-BGR 0x44 0x54 0x6A)  # dk2      44546A
    4  = (New-BGR 0xE7 0xE6 0xE6)  # lt2      E7E6E6
    5  = (New-BGR 0x5B 0x9B 0xD5)  # accent1  5B9BD5
    6  = (New-BGR 0xED 0x7D 0x31)  # accent2  ED7D31
    7  = (New-BGR 0xA5 0xA5 0xA5)  # accent3  A5A5A5
    8  = (New-BGR 0xFF 0xC0 0x00)  # accent4  FFC000
    9  = (New-BGR 0x44 0x72 0xC4)  # accent5  4472C4
    10 = (New-BGR 0x70 0xAD 0x47)  # accent6  70AD47
    11 = (New-BGR 0x05 0x63 0xC1)  # hlink    0563C1
    12 = (New-BGR 0x95 0x4F 0x72)  # folHlink 954F72
}

$slideMasterScheme = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($slot in 1..12) {
    $slideMasterScheme.Colors($slot).RGB = $officeTheme[$slot]
}
